# The commit swaps the two embedded themes:
#   ppt/theme/theme1.xml  (used by the notes master)  "Office Theme" / "Office"   -> "Integral" / "Red Violet"
#   ppt/theme/theme2.xml  (used by the slide master)  "Integral"     / "Red Violet" -> "Office Theme" / "Office"
#
# i.e. the slide master's theme (which is what PowerPoint's object model exposes
# as the deck's single ThemeColorScheme / Design) goes from the pink/violet
# "Integral" palette back to the stock Office palette.
#
# Apply the new ("Office Theme") colours to the presentation's theme colour
# scheme via the documented PowerPoint OM surface (Slide.ThemeColorScheme /
# ThemeColor.RGB), rather than poking the OOXML directly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToComRgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (the palette theme1.xml currently holds).
$tcs.Item(1).RGB  = ToComRgb 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = ToComRgb 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = ToComRgb 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = ToComRgb 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = ToComRgb 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = ToComRgb 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = ToComRgb 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = ToComRgb 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = ToComRgb 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = ToComRgb 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = ToComRgb 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = ToComRgb 0x95 0x4F 0x72   # folHlink
